$d = $word.ActiveDocument

function Set-BoldForText($searchText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
    }
}

# Simple bold-only toggles
Set-BoldForText("Holly Dickson")
Set-BoldForText("애니메이션 인턴")
Set-BoldForText("미술학부 애니메이션과 학사 학위")

# "업무 경력" -> bold on + text shortened to "업무 력"
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("업무 경력", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $range.Font.Bold = 1
    $range.Text = "업무 력"
}

# "보조 애니메이션 디자이너" -> bold on + text changed to "주니어 애니메이션 디자이너".
# This paragraph's mark sits right before the "Graphic Design Institute"
# paragraph; touching the boundary re-stamps that next paragraph with a
# fresh (rsid-less / paraId-less) mark, so we reproduce that by
# deleting + reinserting the paragraph break, then restoring the
# formatting that lived on paragraph 1's mark (outline level + bold).
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("보조 애니메이션 디자이너", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $range.Font.Bold = 1
    $range.Text = "주니어 애니메이션 디자이너"

    $p = $range.Paragraphs(1)
    $origOutlineLevel = $p.OutlineLevel
    $pRange = $p.Range
    $paraStart = $pRange.Start

    $endRange = $pRange.Duplicate
    $endRange.Collapse(0)
    $markOnly = $d.Range($endRange.Start - 1, $endRange.Start)
    $origMarkBold = $markOnly.Font.Bold

    $splitPos = $markOnly.Start
    $markOnly.Delete()
    $insertRange = $d.Range($splitPos, $splitPos)
    $insertRange.InsertParagraphAfter()

    $p1 = $d.Range($paraStart, $paraStart).Paragraphs(1)
    $p1.OutlineLevel = $origOutlineLevel
    $p1Range = $p1.Range
    $p1EndRange = $p1Range.Duplicate
    $p1EndRange.Collapse(0)
    $p1MarkOnly = $d.Range($p1EndRange.Start - 1, $p1EndRange.Start)
    $p1MarkOnly.Font.Bold = $origMarkBold
}
